# Commit: "Fixed update to excel issue"
# - Rename header cell B1 on "Weekly Quantity" from "Requested quantity" to "Weekly_PO_Qty"
# - Rename header cell B1 on "Monthly Trend" from "Requested quantity" to "Monthly_PO_Qty"
# - Add a new worksheet "PO Forecast" (after the existing two sheets) containing a
#   ds / PO_Forecast / yhat_lower / yhat_upper forecast table.

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity sheet: rename header ---------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend sheet: rename header -----------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add new "PO Forecast" sheet at the end -------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @"
45312.99999999999,260,-74.09028283621826,616.7595474675844
45319.99999999999,260,-75.04745294572182,615.2178649297504
45333.99999999999,258,-91.59674654370421,627.6343376509058
45340.99999999999,257,-78.9553813146468,587.0384505591194
45347.99999999999,257,-87.00780581423159,572.6717446991277
45354.99999999999,256,-90.74080596723003,598.2613373277964
45368.99999999999,254,-86.86607770057932,616.4313293992819
45375.99999999999,254,-99.4512911495152,574.6197024585542
45382.99999999999,253,-104.0509961667969,592.3533088636977
45389.99999999999,252,-105.457981448454,602.3025668915893
45396.99999999999,251,-101.4987982507167,592.1878417441991
45403.99999999999,251,-76.80678925544193,605.2905537110588
45417.99999999999,249,-84.85413081937122,592.7264199932772
45424.99999999999,249,-122.1562240361304,606.4185622229065
45431.99999999999,248,-116.649317987699,599.108055333489
45438.99999999999,247,-72.75982988320141,623.0977340235831
45445.99999999999,246,-83.05243639658596,583.4948970784508
45452.99999999999,246,-89.34622342703673,575.5803520328387
45466.99999999999,244,-103.4163053044872,579.3162693663268
45473.99999999999,243,-106.4625370810917,574.6327747404348
45487.99999999999,242,-112.4728905118561,568.1626078663191
45494.99999999999,241,-99.0831913347275,581.9965724644886
45557.99999999999,235,-130.1589827229246,583.1627404163928
45571.99999999999,233,-107.703075432073,572.8578027091633
45578.99999999999,232,-115.2412489228586,570.2600796173292
45592.99999999999,231,-108.5718016576266,569.5754567410589
45599.99999999999,230,-109.3311863916319,580.0064505398418
45606.99999999999,229,-122.3962758980058,600.9495204308572
45613.99999999999,229,-99.50032695944786,571.1479570040548
45620.99999999999,228,-90.9910154087279,557.2510468682502
45627.99999999999,227,-149.7896568311633,573.4823803575476
45634.99999999999,226,-114.0776621243043,588.7749011796632
45641.99999999999,226,-140.2271731815216,559.3917067056352
45648.99999999999,225,-93.59552944340169,562.8021636247668
45655.99999999999,224,-125.2984451269411,574.5164028836153
45662.99999999999,224,-144.5042186546807,544.4485395715994
45669.99999999999,223,-115.7705776649244,563.7517488437637
"@

$lines = $data -split "`n"
$rowIndex = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","

    $dsCell = $wsForecast.Cells.Item($rowIndex, 1)
    $dsCell.Value = [double]$parts[0]
    $dsCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $wsForecast.Cells.Item($rowIndex, 2).Value = [double]$parts[1]
    $wsForecast.Cells.Item($rowIndex, 3).Value = [double]$parts[2]
    $wsForecast.Cells.Item($rowIndex, 4).Value = [double]$parts[3]

    $rowIndex++
}
